$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text values in column B
$ws.Range("B2").Value = "<them>"
$ws.Range("B3").Value = "<his>"

# Update numeric counts in column C
$ws.Range("C6").Value = 6
$ws.Range("C8").Value = 5
$ws.Range("C9").Value = 3
$ws.Range("C11").Value = 4
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 3
$ws.Range("C14").Value = 6
$ws.Range("C15").Value = 2
$ws.Range("C16").Value = 0
$ws.Range("C17").Value = 5
$ws.Range("C18").Value = 2
